# Auto-generated edit script for previsao_retorno.xlsx update
# Commit: "atualizei dados da bibi e add" (updated bibi data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh of active-customer transaction metrics (rows with updated purchase history) ---
# Row 13
$ws.Range("B13").Value = 0.33
$ws.Range("C13").Value = 0.33
$ws.Range("E13").Value = 25
$ws.Range("H13").Value = 45807.49570601852
$ws.Range("I13").Value = 45868.49570601852

# Row 52
$ws.Range("E52").Value = 36
$ws.Range("H52").Value = 45809.63699074074
$ws.Range("I52").Value = 45824.63699074074

# Row 61
$ws.Range("B61").Value = 0.42
$ws.Range("C61").Value = 0.33
$ws.Range("E61").Value = 7
$ws.Range("H61").Value = 45807.40917824074
$ws.Range("I61").Value = 45838.40917824074

# Row 109
$ws.Range("B109").Value = 0.17
$ws.Range("D109").Value = 0.33
$ws.Range("E109").Value = 13
$ws.Range("F109").Value = 0.33
$ws.Range("H109").Value = 45807.45422453704
$ws.Range("I109").Value = 45868.45422453704

# Row 111
$ws.Range("E111").Value = 15038
$ws.Range("H111").Value = 45807.76511574074
$ws.Range("I111").Value = 45808.76511574074

# --- Refresh of "months since last purchase" labels for inactive customers ---
$situacaoUpdates = @{
    2 = "INATIVO - 54.4 meses sem comprar"
    4 = "INATIVO - 35.2 meses sem comprar"
    5 = "INATIVO - 13.9 meses sem comprar"
    6 = "INATIVO - 15.4 meses sem comprar"
    8 = "INATIVO - 17.1 meses sem comprar"
    9 = "INATIVO - 18.2 meses sem comprar"
    10 = "INATIVO - 1.7 meses sem comprar"
    11 = "INATIVO - 4.1 meses sem comprar"
    15 = "INATIVO - 38.9 meses sem comprar"
    16 = "INATIVO - 1.6 meses sem comprar"
    17 = "INATIVO - 36.5 meses sem comprar"
    18 = "INATIVO - 10.5 meses sem comprar"
    19 = "INATIVO - 14.4 meses sem comprar"
    20 = "INATIVO - 36.4 meses sem comprar"
    22 = "INATIVO - 28.3 meses sem comprar"
    23 = "INATIVO - 36.9 meses sem comprar"
    25 = "INATIVO - 18.6 meses sem comprar"
    27 = "INATIVO - 16.7 meses sem comprar"
    30 = "INATIVO - 21.6 meses sem comprar"
    31 = "INATIVO - 6.7 meses sem comprar"
    33 = "INATIVO - 13.7 meses sem comprar"
    34 = "INATIVO - 25.9 meses sem comprar"
    37 = "INATIVO - 31.3 meses sem comprar"
    38 = "INATIVO - 33.5 meses sem comprar"
    39 = "INATIVO - 11.6 meses sem comprar"
    41 = "INATIVO - 6.6 meses sem comprar"
    42 = "INATIVO - 24.9 meses sem comprar"
    44 = "INATIVO - 15.0 meses sem comprar"
    45 = "INATIVO - 1.6 meses sem comprar"
    48 = "INATIVO - 6.6 meses sem comprar"
    49 = "INATIVO - 9.1 meses sem comprar"
    57 = "INATIVO - 11.0 meses sem comprar"
    63 = "INATIVO - 26.9 meses sem comprar"
    64 = "INATIVO - 21.0 meses sem comprar"
    66 = "INATIVO - 11.7 meses sem comprar"
    67 = "INATIVO - 10.5 meses sem comprar"
    71 = "INATIVO - 0.4 meses sem comprar"
    72 = "INATIVO - 21.0 meses sem comprar"
    73 = "INATIVO - 32.5 meses sem comprar"
    74 = "INATIVO - 6.8 meses sem comprar"
    79 = "INATIVO - 21.4 meses sem comprar"
    81 = "INATIVO - 25.3 meses sem comprar"
    83 = "INATIVO - 20.9 meses sem comprar"
    84 = "INATIVO - 8.6 meses sem comprar"
    85 = "INATIVO - 14.1 meses sem comprar"
    86 = "INATIVO - 4.0 meses sem comprar"
    87 = "INATIVO - 10.9 meses sem comprar"
    88 = "INATIVO - 10.3 meses sem comprar"
    89 = "INATIVO - 14.2 meses sem comprar"
    90 = "INATIVO - 32.5 meses sem comprar"
    91 = "INATIVO - 12.8 meses sem comprar"
    92 = "INATIVO - 17.8 meses sem comprar"
    93 = "INATIVO - 15.5 meses sem comprar"
    94 = "INATIVO - 18.3 meses sem comprar"
    95 = "INATIVO - 32.0 meses sem comprar"
    97 = "INATIVO - 1.3 meses sem comprar"
    98 = "INATIVO - 22.0 meses sem comprar"
    99 = "INATIVO - 36.3 meses sem comprar"
    100 = "INATIVO - 6.6 meses sem comprar"
    101 = "INATIVO - 13.6 meses sem comprar"
    102 = "INATIVO - 24.0 meses sem comprar"
    103 = "INATIVO - 9.6 meses sem comprar"
    104 = "INATIVO - 24.4 meses sem comprar"
    105 = "INATIVO - 14.8 meses sem comprar"
    106 = "INATIVO - 5.7 meses sem comprar"
    107 = "INATIVO - 20.5 meses sem comprar"
    108 = "INATIVO - 6.2 meses sem comprar"
    110 = "INATIVO - 7.6 meses sem comprar"
}
foreach ($row in $situacaoUpdates.Keys) {
    $ws.Range("J$row").Value = $situacaoUpdates[$row]
}
